$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3708.5535
$ws.Range("I64").Value = 3550
$ws.Range("J64").Value = 3845.9666
$ws.Range("K64").Value = 3550
$ws.Range("L64").Value = 3845.9666
$ws.Range("M64").Value = -3302
$ws.Range("N64").Value = -4341.9666
$ws.Range("H67").Value = 3708.5535
$ws.Range("I67").Value = 3550
$ws.Range("J67").Value = 3845.9666
$ws.Range("K67").Value = 3550
$ws.Range("L67").Value = 3845.9666
$ws.Range("M67").Value = -2692
$ws.Range("N67").Value = -5561.9666
$ws.Range("H98").Value = 2232.375
$ws.Range("I98").Value = 1781
$ws.Range("J98").Value = 2984.6667
$ws.Range("K98").Value = 1781
$ws.Range("L98").Value = 2984.6667
$ws.Range("M98").Value = -283
$ws.Range("N98").Value = -5980.6667
$ws.Range("H106").Value = 8746.75
$ws.Range("I106").Value = 4990.5
$ws.Range("J106").Value = 12503
$ws.Range("K106").Value = 4990.5
$ws.Range("L106").Value = 12503
$ws.Range("M106").Value = -4359.5
$ws.Range("N106").Value = -13765
$ws.Range("H122").Value = 2232.375
$ws.Range("I122").Value = 1781
$ws.Range("J122").Value = 2984.6667
$ws.Range("K122").Value = 5343
$ws.Range("L122").Value = 8954.000100000001
$ws.Range("M122").Value = -2893
$ws.Range("N122").Value = -13854.0001
$ws.Range("H137").Value = 1573.8948
$ws.Range("I137").Value = 880.9032
$ws.Range("J137").Value = 4642.857
$ws.Range("K137").Value = 2642.7096
$ws.Range("L137").Value = 13928.571
$ws.Range("M137").Value = -92.70960000000014
$ws.Range("N137").Value = -19028.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 240.11111
$ws.Range("I5").Value = 265.57144
$ws.Range("K5").Value = 265.57144
$ws.Range("M5").Value = -153.57144
$ws.Range("H32").Value = 16963300
$ws.Range("I32").Value = 23260928
$ws.Range("J32").Value = 38425.812
$ws.Range("K32").Value = 23260928
$ws.Range("L32").Value = 38425.812
$ws.Range("M32").Value = -23260641
$ws.Range("N32").Value = -38999.812
$ws.Range("H45").Value = 1364.2354
$ws.Range("I45").Value = 977
$ws.Range("J45").Value = 1917.4286
$ws.Range("K45").Value = 977
$ws.Range("L45").Value = 1917.4286
$ws.Range("M45").Value = -600
$ws.Range("N45").Value = -2671.4286
$ws.Range("H125").Value = 58750
$ws.Range("J125").Value = 58750
$ws.Range("L125").Value = 58750
$ws.Range("N125").Value = -68590
$ws.Range("H132").Value = 1793.7931
$ws.Range("I132").Value = 1117.4445
$ws.Range("J132").Value = 2900.5454
$ws.Range("K132").Value = 3352.3335
$ws.Range("L132").Value = 8701.636200000001
$ws.Range("M132").Value = -822.3335000000002
$ws.Range("N132").Value = -13761.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 240.11111
$ws.Range("I4").Value = 265.57144
$ws.Range("K4").Value = 265.57144
$ws.Range("M4").Value = -150.57144
$ws.Range("H134").Value = 1635.6216
$ws.Range("I134").Value = 1374.6666
$ws.Range("J134").Value = 2117.3845
$ws.Range("K134").Value = 4123.9998
$ws.Range("L134").Value = 6352.1535
$ws.Range("M134").Value = -1588.9998
$ws.Range("N134").Value = -11422.1535

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3323.9375
$ws.Range("I31").Value = 2304
$ws.Range("J31").Value = 4635.2856
$ws.Range("K31").Value = 2304
$ws.Range("L31").Value = 4635.2856
$ws.Range("M31").Value = -2009
$ws.Range("N31").Value = -5225.2856
$ws.Range("H34").Value = 3323.9375
$ws.Range("I34").Value = 2304
$ws.Range("J34").Value = 4635.2856
$ws.Range("K34").Value = 2304
$ws.Range("L34").Value = 4635.2856
$ws.Range("M34").Value = -2102
$ws.Range("N34").Value = -5039.2856
$ws.Range("H62").Value = 3037.8462
$ws.Range("I62").Value = 2680
$ws.Range("J62").Value = 3344.5715
$ws.Range("K62").Value = 2680
$ws.Range("L62").Value = 3344.5715
$ws.Range("M62").Value = -2056
$ws.Range("N62").Value = -4592.5715
$ws.Range("H65").Value = 3037.8462
$ws.Range("I65").Value = 2680
$ws.Range("J65").Value = 3344.5715
$ws.Range("K65").Value = 13400
$ws.Range("L65").Value = 16722.8575
$ws.Range("M65").Value = -10280
$ws.Range("N65").Value = -22962.8575
$ws.Range("H122").Value = 1647.3636
$ws.Range("I122").Value = 1058.8
$ws.Range("J122").Value = 2908.5715
$ws.Range("K122").Value = 3176.4
$ws.Range("L122").Value = 8725.7145
$ws.Range("M122").Value = -726.3999999999996
$ws.Range("N122").Value = -13625.7145
$ws.Range("H134").Value = 2040.8182
$ws.Range("I134").Value = 1327.0227
$ws.Range("K134").Value = 3981.0681
$ws.Range("M134").Value = -1446.0681

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 6.869565
$ws.Range("I12").Value = 17.833334
$ws.Range("J12").Value = 3
$ws.Range("K12").Value = 53.500002
$ws.Range("L12").Value = 9
$ws.Range("M12").Value = 119.499998
$ws.Range("N12").Value = -355

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1650.5312
$ws.Range("I102").Value = 1546.8077
$ws.Range("J102").Value = 2100
$ws.Range("K102").Value = 1546.8077
$ws.Range("L102").Value = 2100
$ws.Range("M102").Value = 75.19229999999993
$ws.Range("N102").Value = -5344
$ws.Range("H132").Value = 2919.5193
$ws.Range("I132").Value = 2625.1794
$ws.Range("J132").Value = 3802.5386
$ws.Range("K132").Value = 7875.5382
$ws.Range("L132").Value = 11407.6158
$ws.Range("M132").Value = -5345.5382
$ws.Range("N132").Value = -16467.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 574.875
$ws.Range("I46").Value = 100
$ws.Range("J46").Value = 1049.75
$ws.Range("K46").Value = 100
$ws.Range("L46").Value = 1049.75
$ws.Range("M46").Value = 88
$ws.Range("N46").Value = -1425.75
$ws.Range("H55").Value = 227.53847
$ws.Range("I55").Value = 132.81818
$ws.Range("J55").Value = 297
$ws.Range("K55").Value = 132.81818
$ws.Range("L55").Value = 297
$ws.Range("M55").Value = 40.18181999999999
$ws.Range("N55").Value = -643
$ws.Range("H82").Value = 2822.9092
$ws.Range("I82").Value = 1734
$ws.Range("J82").Value = 3231.25
$ws.Range("K82").Value = 1734
$ws.Range("L82").Value = 3231.25
$ws.Range("M82").Value = -1373
$ws.Range("N82").Value = -3953.25
$ws.Range("H85").Value = 2822.9092
$ws.Range("I85").Value = 1734
$ws.Range("J85").Value = 3231.25
$ws.Range("K85").Value = 1734
$ws.Range("L85").Value = 3231.25
$ws.Range("M85").Value = -486
$ws.Range("N85").Value = -5727.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1988.3125
$ws.Range("I126").Value = 1539.2222
$ws.Range("J126").Value = 2565.7144
$ws.Range("K126").Value = 4617.6666
$ws.Range("L126").Value = 7697.1432
$ws.Range("M126").Value = -2147.6666
$ws.Range("N126").Value = -12637.1432
$ws.Range("H132").Value = 2094.5476
$ws.Range("I132").Value = 1048.2693
$ws.Range("J132").Value = 3794.75
$ws.Range("K132").Value = 3144.8079
$ws.Range("L132").Value = 11384.25
$ws.Range("M132").Value = -614.8078999999998
$ws.Range("N132").Value = -16444.25

Write-Host "Done applying Shinryu_Profits updates"